# Applies the "ifo GDP component analysis preprocessing" update:
# extends the staircase error matrix by one more diagonal column per row
# (rows 7-16), recomputing the previous edge cell and adding the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: recompute existing last value (K6)
$ws.Range("K6").Value = -0.2254024683979639

# Row 7: recompute previous edge (J7) and add new edge (K7)
$ws.Range("J7").Value = -0.1253231084953424
$ws.Range("K7").Value = -0.3352267436446591

# Row 8: recompute previous edge (I8) and add new edge (J8)
$ws.Range("I8").Value = 0.2284633975843539
$ws.Range("J8").Value = 0.01855976243503714

# Row 9: recompute previous edge (H9) and add new edge (I9)
$ws.Range("H9").Value = 0.08028600715190851
$ws.Range("I9").Value = -0.1296176279974082

# Row 10: recompute previous edge (G10) and add new edge (H10)
$ws.Range("G10").Value = -0.07715998185224648
$ws.Range("H10").Value = -0.2870636170015632

# Row 11: recompute previous edge (F11) and add new edge (G11)
$ws.Range("F11").Value = 0.4234994746738243
$ws.Range("G11").Value = 0.2135958395245076

# Row 12: recompute previous edge (E12) and add new edge (F12)
$ws.Range("E12").Value = 0.1431415941383551
$ws.Range("F12").Value = -0.06676204101096155

# Row 13: recompute previous edge (D13) and add new edge (E13)
$ws.Range("D13").Value = 0.3151164519833668
$ws.Range("E13").Value = 0.1052128168340501

# Row 14: recompute previous edge (C14) and add new edge (D14)
$ws.Range("C14").Value = 0.009253912237035311
$ws.Range("D14").Value = -0.2006497229122814

# Row 15: recompute previous edge (B15) and add new edge (C15)
$ws.Range("B15").Value = 0.6215838649243215
$ws.Range("C15").Value = 0.4116802297750048

# Row 16: add new edge (B16)
$ws.Range("B16").Value = -0.2766911554241067
